# Fruta / hortaliza, semanal
# Insert a new weekly record at row 100 (Feria Lagunitas de Puerto Montt - Piña),
# pushing the existing rows 100-144 down to 101-145.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 100; Excel shifts row 100..144 -> 101..145
# and copies formatting (e.g. the date style on column D) down from the row above.
$ws.Rows.Item(100).Insert()

$ws.Range("A100").Value = 4
$ws.Range("B100").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C100").Value = "Los Lagos"
$ws.Range("D100").Value = 44489
$ws.Range("E100").Value = 10
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100108
$ws.Range("H100").Value = "Tropicales y subtropicales"
$ws.Range("I100").Value = 100108005
$ws.Range("J100").Value = "Piña"
$ws.Range("K100").Value = "Caramelo"
$ws.Range("L100").Value = "Segunda"
$ws.Range("M100").Value = 30
$ws.Range("N100").Value = 25000
$ws.Range("O100").Value = 25000
$ws.Range("P100").Value = 25000
$ws.Range("Q100").Value = "`$/caja 14 unidades"
$ws.Range("R100").Value = "Ecuador"
$ws.Range("S100").Value = 1786
$ws.Range("T100").Value = 14
